# "Add files via upload" - refresh the W9 Salaries and Tasks sheet with this
# week's data: date, team name, member roster + salary split, and the
# completed/upcoming task lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header block -----------------------------------------------------
$ws.Range("B3").Value = 43776  # Date: 11/7/2019
$ws.Range("B5").Value = 5       # Total Number of Team Members

# --- Team member names + individual salary -----------------------------
$ws.Range("A8").Value = "Kunaal Sikka"
$ws.Range("A9").Value = "Mina Huh"
$ws.Range("A10").Value = "Vu Nguyen"
$ws.Range("A11").Value = "Nicolas Carmody"
$ws.Range("A12").Value = "Jonas Bokstaller"

$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 100

$ws.Range("B4").Value = "MSR Voice Input" # Team Name

# --- Tasks completed this week / tasks to complete next week -----------
$ws.Range("A19").Value = "Create Questionnaire for experiment"
$ws.Range("A20").Value = "Design experiment VA/Documentation"
$ws.Range("A21").Value = "Test experiment with members "

$ws.Range("B19").Value = "Perform experiment on other groups"
$ws.Range("B20").Value = "Evaluate experiment/questionnaire"

# --- Misc cosmetic/view updates matching the author's session ----------
$wb.Styles.Item(1).Name = "Standard"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("B20").Select()
